# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r, B(TB), C(d2S), D(K), E(IP)  -- F(Win) unchanged, G(sum) = B+C+D+E
$data = @(
    @(2, 3.182878228561681, 86.29678392075563, 0.1529057820181812, 6.48142807727062),
    @(3, 0.7287194209349384, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569),
    @(4, 0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569),
    @(5, 1.505614041169197, 86.29678392075563, 0.7127328510149897, 246.9852506941017),
    @(6, 0.1554434735375247, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569),
    @(7, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569),
    @(8, 3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569),
    @(9, 1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569),
    @(10, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569),
    @(11, 0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569),
    @(12, 0.02258322285507441, 0.05231270169004087, 3.082599426703578, 0.4998867070740569),
    @(13, 3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569),
    @(14, 0.3464964993005633, 0.05231270169004087, 16.98373111632243, 6.48142807727062),
    @(15, 0.06328177979961902, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569),
    @(16, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569),
    @(17, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569),
    @(18, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569),
    @(19, 0.06328177979961902, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569),
    @(20, 1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569),
    @(21, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569),
    @(22, 1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569),
    @(23, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569),
    @(24, 3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569),
    @(25, 0.3464964993005633, 1.65323645889881, 0.7127328510149897, 0.4998867070740569)
)

foreach ($row in $data) {
    $r = $row[0]
    $b = $row[1]
    $c = $row[2]
    $d = $row[3]
    $e = $row[4]

    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 7).Value = $b + $c + $d + $e
}
